# Updates team-specific time matrix values (Georgetown_A) on Sheet1.
# Each assignment below sets a single data cell to its new decimal value
# as described by the source diff (row labels in column A are unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1869158878504673
$ws.Range("C2").Value = 0.5607476635514018
$ws.Range("J2").Value = 0.006230529595015576
$ws.Range("P2").Value = 0.1370716510903427
$ws.Range("S2").Value = 0.1090342679127726

# Row 3
$ws.Range("B3").Value = 0.02116402116402116
$ws.Range("C3").Value = 0.04761904761904762
$ws.Range("J3").Value = 0.01587301587301587
$ws.Range("P3").Value = 0.7037037037037037
$ws.Range("S3").Value = 0.2116402116402116

# Row 4
$ws.Range("P4").Value = 0.8857142857142857
$ws.Range("S4").Value = 0.1142857142857143

# Row 6
$ws.Range("B6").Value = 0.0430622009569378
$ws.Range("D6").Value = 0.004784688995215311
$ws.Range("E6").Value = 0.009569377990430622
$ws.Range("F6").Value = 0.05263157894736842
$ws.Range("J6").Value = 0.2392344497607655
$ws.Range("O6").Value = 0.03827751196172249
$ws.Range("Q6").Value = 0.1531100478468899
$ws.Range("R6").Value = 0.09090909090909091
$ws.Range("S6").Value = 0.3684210526315789

# Row 7
$ws.Range("B7").Value = 0.1288659793814433
$ws.Range("D7").Value = 0.01030927835051546
$ws.Range("E7").Value = 0.005154639175257732
$ws.Range("F7").Value = 0.07731958762886598
$ws.Range("J7").Value = 0.1649484536082474
$ws.Range("O7").Value = 0.0154639175257732
$ws.Range("Q7").Value = 0.1597938144329897
$ws.Range("R7").Value = 0.1082474226804124
$ws.Range("S7").Value = 0.3298969072164948

# Row 8
$ws.Range("B8").Value = 0.1067961165048544
$ws.Range("D8").Value = 0.01359223300970874
$ws.Range("F8").Value = 0.05825242718446602
$ws.Range("J8").Value = 0.1378640776699029
$ws.Range("O8").Value = 0.01747572815533981
$ws.Range("Q8").Value = 0.1495145631067961
$ws.Range("R8").Value = 0.0912621359223301
$ws.Range("S8").Value = 0.4252427184466019

# Row 9
$ws.Range("B9").Value = 0.1244239631336406
$ws.Range("D9").Value = 0.03686635944700461
$ws.Range("F9").Value = 0.06451612903225806
$ws.Range("J9").Value = 0.119815668202765
$ws.Range("O9").Value = 0.01382488479262673
$ws.Range("Q9").Value = 0.1566820276497696
$ws.Range("R9").Value = 0.06912442396313365
$ws.Range("S9").Value = 0.4147465437788018

# Row 10
$ws.Range("B10").Value = 0.1142618849040867
$ws.Range("D10").Value = 0.0158465387823186
$ws.Range("E10").Value = 0.00250208507089241
$ws.Range("F10").Value = 0.07172643869891576
$ws.Range("J10").Value = 0.1359466221851543
$ws.Range("O10").Value = 0.01751459549624687
$ws.Range("Q10").Value = 0.1843202668890742
$ws.Range("R10").Value = 0.09257714762301918
$ws.Range("S10").Value = 0.3653044203502919

# Row 11
$ws.Range("G11").Value = 0.1496598639455782
$ws.Range("J11").Value = 0.1054421768707483
$ws.Range("K11").Value = 0.2312925170068027
$ws.Range("L11").Value = 0.5
$ws.Range("S11").Value = 0.01360544217687075

# Row 12
$ws.Range("G12").Value = 0.7453416149068323
$ws.Range("J12").Value = 0.1614906832298137
$ws.Range("K12").Value = 0.006211180124223602
$ws.Range("L12").Value = 0.03726708074534162
$ws.Range("S12").Value = 0.04968944099378882

# Row 13
$ws.Range("G13").Value = 0.7608695652173914
$ws.Range("J13").Value = 0.1521739130434783
$ws.Range("S13").Value = 0.08695652173913043

# Row 15
$ws.Range("F15").Value = 0.03211009174311927
$ws.Range("H15").Value = 0.1422018348623853
$ws.Range("I15").Value = 0.0871559633027523
$ws.Range("J15").Value = 0.3532110091743119
$ws.Range("K15").Value = 0.05045871559633028
$ws.Range("M15").Value = 0.01834862385321101
$ws.Range("O15").Value = 0.04128440366972477
$ws.Range("S15").Value = 0.2752293577981652

# Row 16
$ws.Range("F16").Value = 0.005050505050505051
$ws.Range("H16").Value = 0.2424242424242424
$ws.Range("I16").Value = 0.1161616161616162
$ws.Range("J16").Value = 0.3181818181818182
$ws.Range("K16").Value = 0.1111111111111111
$ws.Range("M16").Value = 0.0303030303030303
$ws.Range("O16").Value = 0.04545454545454546
$ws.Range("S16").Value = 0.1313131313131313

# Row 17
$ws.Range("F17").Value = 0.0100250626566416
$ws.Range("H17").Value = 0.2155388471177945
$ws.Range("I17").Value = 0.09523809523809523
$ws.Range("J17").Value = 0.3809523809523809
$ws.Range("K17").Value = 0.1228070175438596
$ws.Range("M17").Value = 0.02255639097744361
$ws.Range("O17").Value = 0.05012531328320802
$ws.Range("S17").Value = 0.1027568922305764

# Row 18
$ws.Range("F18").Value = 0.01388888888888889
$ws.Range("H18").Value = 0.2175925925925926
$ws.Range("I18").Value = 0.09259259259259259
$ws.Range("J18").Value = 0.3518518518518519
$ws.Range("K18").Value = 0.1064814814814815
$ws.Range("M18").Value = 0.009259259259259259
$ws.Range("O18").Value = 0.06481481481481481
$ws.Range("S18").Value = 0.1435185185185185

# Row 19
$ws.Range("F19").Value = 0.01157407407407407
$ws.Range("H19").Value = 0.2376543209876543
$ws.Range("I19").Value = 0.08950617283950617
$ws.Range("J19").Value = 0.3371913580246914
$ws.Range("K19").Value = 0.08719135802469136
$ws.Range("M19").Value = 0.01929012345679012
$ws.Range("O19").Value = 0.07561728395061729
$ws.Range("S19").Value = 0.1419753086419753
